$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing existing rows 3-10 down to 4-11
$ws.Rows("3:3").Insert()

# Mark row 1's item ("A la place du champ description...") as done
$ws.Range("C1").Value = "OK"

# The "OK" that used to live on row 2 no longer applies there
$ws.Range("C2").ClearContents()

# New row 3: the new "use video for profile photo" feature, marked done
$ws.Range("B3").Value = "Afficher une vidéo de profil au lieu d'une photo"
$ws.Range("C3").Value = "OK"

# Update the active selection to match the saved view state
$ws.Range("C4").Select()
